$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.564.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.18%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.32%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.76"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.462.40"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.985.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.612.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.59"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.22"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.84"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.15"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.36"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.51"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0443"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.95%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.01%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.60%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.272"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.041.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.535"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +15.96%  "
